$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 42
$link = "https://www.360dx.com/diagnostics/proscia-using-50m-fundraise-expand-staff-bolster-ai-tools-menu"
$keywords = "digital pathology"
$title = "Proscia Using `$50M Fundraise to Expand Staff, Bolster AI Tools Menu"

$ws.Cells.Item($newRow, 1).Value = $link
$ws.Cells.Item($newRow, 2).Value = $keywords
$ws.Cells.Item($newRow, 3).Value = $title

$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $link) | Out-Null
$ws.Cells.Item($newRow, 1).Style = "Hyperlink"
